# Regenerate save_data column G ("K") values: replace the previous
# "Strike#" derived values with newly calculated K values (std/mean based
# s_vals calculation regenerated these numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    3, 4, 3, 4, 4, 6, 7, 2, 2, 4,
    4, 6, 2, 8, 6, 5, 8, 5, 8, 2,
    1, 3, 2, 3, 6, 6, 6, 3, 4, 4,
    5, 5, 5, 2
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
